$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update credentials for the "admin" (non-test) user/connection string.
$ws.Range("B9").Value = "105881_av56092"
$ws.Range("B10").Value = "105881_mj42632"

$ws.Range("B21").Value = "ODBC;DRIVER=SQL Server Native Client 10.0;SERVER=admin-105881.mssql.stwcp.net;UID=105881_kr28513;PWD=Maine1953;APP=Microsoft Office;DATABASE=10588_admin;"
$ws.Range("B22").Value = "ODBC;DRIVER=SQL Server Native Client 10.0;SERVER=admintest-105881.mssql.stwcp.net;UID=105881-mj42632;PWD=Maine1953;APP=Microsoft Office;DATABASE=105881-admintest;"

# Update view state: selection moved from B4 to B10 (sheet also scrolled so
# row 4 is now the top-left visible row).
$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("B10").Select()
